$wb = $excel.ActiveWorkbook

$wsPerf = $wb.Worksheets.Item("Perf")
$wsClass = $wb.Worksheets.Item("Class")

# Rename the left/right split labels on both sheets:
#   "Left target"  -> "Left met"
#   "Left met"     -> "Left not met"
#   "Right target" -> "Right met"
#   "Right met"    -> "Right not met"
foreach ($ws in @($wsPerf, $wsClass)) {
    $ws.Range("E4").Value = "Left met"
    $ws.Range("E5").Value = "Left not met"
    $ws.Range("E7").Value = "Right met"
    $ws.Range("E8").Value = "Right not met"
}

# Switch the active tab/selection from "Class" to "Perf", both now pointing at E9
$wsClass.Range("E9").Select()
$wsPerf.Activate()
$wsPerf.Range("E9").Select()
